$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 863, shifting rows 863:935 down to 864:936
$ws.Rows.Item(863).EntireRow.Insert()

# Populate the newly inserted row 863 (copy of neighboring row's fixed columns
# plus the new record's own Fecha / Variedad / Volumen)
$ws.Cells.Item(863, 1).Value = 10
$ws.Cells.Item(863, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(863, 3).Value = "La Araucanía"
$ws.Cells.Item(863, 4).Value = 44783
$ws.Cells.Item(863, 5).Value = 9
$ws.Cells.Item(863, 6).Value = "Fruta"
$ws.Cells.Item(863, 7).Value = 100102
$ws.Cells.Item(863, 8).Value = "Cítricos"
$ws.Cells.Item(863, 9).Value = 100102005
$ws.Cells.Item(863, 10).Value = "Naranja"
$ws.Cells.Item(863, 11).Value = "New Hall"
$ws.Cells.Item(863, 12).Value = "Primera"
$ws.Cells.Item(863, 13).Value = 215
$ws.Cells.Item(863, 14).Value = 9000
$ws.Cells.Item(863, 15).Value = 9000
$ws.Cells.Item(863, 16).Value = 9000
$ws.Cells.Item(863, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(863, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(863, 19).Value = 600
$ws.Cells.Item(863, 20).Value = 15
